$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 6) - data now only spans rows 1-5
$ws.Rows.Item(6).Delete()

# Overwrite rows 2-5 with the refreshed sensor readings (new dataset slice)
$ws.Range("A2").Value = 45100.50694444445
$ws.Range("B2").Value = 22.1
$ws.Range("C2").Value = 15.162
$ws.Range("D2").Value = 4.209
$ws.Range("E2").Value = 46.427
$ws.Range("F2").Value = 38.398
$ws.Range("G2").Value = 17.391
$ws.Range("H2").Value = 57.345
$ws.Range("I2").Value = 26.759
$ws.Range("J2").Value = 11.337
$ws.Range("K2").Value = 17.452
$ws.Range("L2").Value = 18.439
$ws.Range("M2").Value = 19.302
$ws.Range("N2").Value = 5.553
$ws.Range("O2").Value = 17.294
$ws.Range("P2").Value = 24.297
$ws.Range("Q2").Value = 14.517
$ws.Range("R2").Value = 3.762
$ws.Range("S2").Value = 2.458
$ws.Range("T2").Value = 255.893
$ws.Range("U2").Value = 48.146
$ws.Range("V2").Value = 15.963
$ws.Range("W2").Value = 31.894
$ws.Range("X2").Value = 16.657
$ws.Range("Y2").Value = 2.022
$ws.Range("Z2").Value = 28.586
$ws.Range("AA2").Value = 14.1
$ws.Range("AB2").Value = 12.689
$ws.Range("AC2").Value = 14.848
$ws.Range("AD2").Value = 19.074
$ws.Range("AE2").Value = 3.64
$ws.Range("AF2").Value = 50.56
$ws.Range("AG2").Value = 8.856
$ws.Range("AH2").Value = 19.957

$ws.Range("A3").Value = 45100.51388888889
$ws.Range("B3").Value = 13.452
$ws.Range("C3").Value = 9.414
$ws.Range("D3").Value = 1.718
$ws.Range("E3").Value = 28.608
$ws.Range("F3").Value = 23.668
$ws.Range("G3").Value = 10.586
$ws.Range("H3").Value = 42.895
$ws.Range("I3").Value = 16.288
$ws.Range("J3").Value = 7.001
$ws.Range("K3").Value = 10.581
$ws.Range("L3").Value = 11.535
$ws.Range("M3").Value = 12.001
$ws.Range("N3").Value = 3.383
$ws.Range("O3").Value = 10.527
$ws.Range("P3").Value = 14.828
$ws.Range("Q3").Value = 9.055
$ws.Range("R3").Value = 1.594
$ws.Range("S3").Value = 0.978
$ws.Range("T3").Value = 152.919
$ws.Range("U3").Value = 29.521
$ws.Range("V3").Value = 9.717
$ws.Range("W3").Value = 19.527
$ws.Range("X3").Value = 10.455
$ws.Range("Y3").Value = 1.201
$ws.Range("Z3").Value = 20.432
$ws.Range("AA3").Value = 8.583
$ws.Range("AB3").Value = 7.795
$ws.Range("AC3").Value = 9.122
$ws.Range("AD3").Value = 11.993
$ws.Range("AE3").Value = 1.294
$ws.Range("AF3").Value = 38.79
$ws.Range("AG3").Value = 5.361
$ws.Range("AH3").Value = 12.148

$ws.Range("A4").Value = 45100.52083333334
$ws.Range("B4").Value = 14.893
$ws.Range("C4").Value = 10.734
$ws.Range("D4").Value = 1.273
$ws.Range("E4").Value = 31.973
$ws.Range("F4").Value = 26.444
$ws.Range("G4").Value = 11.72
$ws.Range("H4").Value = 45.754
$ws.Range("I4").Value = 18.033
$ws.Range("J4").Value = 7.887
$ws.Range("K4").Value = 11.831
$ws.Range("L4").Value = 12.926
$ws.Range("M4").Value = 13.493
$ws.Range("N4").Value = 3.744
$ws.Range("O4").Value = 11.655
$ws.Range("P4").Value = 16.495
$ws.Range("Q4").Value = 9.94
$ws.Range("R4").Value = 1.068
$ws.Range("S4").Value = 0.769
$ws.Range("T4").Value = 170.073
$ws.Range("U4").Value = 32.594
$ws.Range("V4").Value = 10.758
$ws.Range("W4").Value = 21.735
$ws.Range("X4").Value = 11.635
$ws.Range("Y4").Value = 1.387
$ws.Range("Z4").Value = 21.991
$ws.Range("AA4").Value = 9.502
$ws.Range("AB4").Value = 8.529
$ws.Range("AC4").Value = 10.003
$ws.Range("AD4").Value = 13.496
$ws.Range("AE4").Value = 0.784
$ws.Range("AF4").Value = 41.255
$ws.Range("AG4").Value = 5.994
$ws.Range("AH4").Value = 13.45

$ws.Range("A5").Value = 45100.52777777778
$ws.Range("B5").Value = 1.92
$ws.Range("C5").Value = 1.11
$ws.Range("D5").Value = 0.62
$ws.Range("E5").Value = 3.86
$ws.Range("F5").Value = 3.2
$ws.Range("G5").Value = 1.52
$ws.Range("H5").Value = 12.4
$ws.Range("I5").Value = 2.33
$ws.Range("J5").Value = 1.01
$ws.Range("K5").Value = 1.34
$ws.Range("L5").Value = 1.65
$ws.Range("M5").Value = 1.58
$ws.Range("N5").Value = 0.51
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.18
$ws.Range("Q5").Value = 1.49
$ws.Range("R5").Value = 0.71
$ws.Range("S5").Value = 0.28
$ws.Range("T5").Value = 15.68
$ws.Range("U5").Value = 4.66
$ws.Range("V5").Value = 1.39
$ws.Range("W5").Value = 3.04
$ws.Range("X5").Value = 1.67
$ws.Range("Y5").Value = 0.08
$ws.Range("Z5").Value = 5.4
$ws.Range("AA5").Value = 1.23
$ws.Range("AB5").Value = 1.24
$ws.Range("AC5").Value = 1.42
$ws.Range("AD5").Value = 1.66
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 11.61
$ws.Range("AG5").Value = 0.67
$ws.Range("AH5").Value = 1.74

# Column width tweaks that came along with the refreshed data (custom accuracy pass)
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667  # B: 7 -> 8
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667  # C: 7 -> 8
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667  # F: 7 -> 8
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667  # G: 7 -> 8
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667  # I: 7 -> 8
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667  # J: 7 -> 8
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667  # K: 7 -> 8
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667  # L: 7 -> 8
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667  # M: 7 -> 8
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667  # O: 7 -> 8
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667  # P: 7 -> 8
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667  # Q: 7 -> 8
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666  # T: 8 -> 9
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667  # V: 7 -> 8
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667  # W: 7 -> 8
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667  # X: 7 -> 8
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667  # Z: 7 -> 8
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667  # AB: 7 -> 8
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667  # AC: 7 -> 8
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667  # AD: 7 -> 8
$ws.Columns.Item(32).ColumnWidth = 7.166666666666667  # AF: 7 -> 8
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667  # AH: 7 -> 8
